# Add a new column R (year 2021) to the worksheet, mirroring the existing
# column Q (year 2020) for formatting, then fill in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column Q (rows 2-13) into column R so the new
# column inherits the same borders/number formats/fonts as the rest of
# the table.
$ws.Range("Q2:Q13").Copy()
$ws.Range("R2:R13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Row 3 - year header
$ws.Range("R3").Value = 2021

# Row 4 - headline indicator (uses a dedicated "0.0" number format, bold
# font, like the rest of row 4)
$ws.Range("R4").NumberFormat = "0.0"
$ws.Range("R4").Value = 18

# Rows 5-13 - regional values for 2021
$ws.Range("R5").Value = 1.7480265877296817
$ws.Range("R6").Value = 4.1112601249414027
$ws.Range("R7").Value = 1.5225742120245318
$ws.Range("R8").Value = 1.2326518235454269
$ws.Range("R9").Value = 4.0865392096984241
$ws.Range("R10").Value = 0.84876624403485645
$ws.Range("R11").Value = 2.1456657699653627
$ws.Range("R12").Value = 1.8214779402142154
$ws.Range("R13").Value = 0.51989507542472779

# Update the selection to mirror the authored workbook (cursor moved to
# R24:R25 after the edits were made).
[void]$ws.Range("R24:R25").Select()
